$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 4.58 = 18113.55 pesos`n✅ 18113.55 pesos = 4.56 = 958.81 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Sheet "tasas": update rate values ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 218.4
$ws2.Range("N12").Value = 3972
$ws2.Range("O12").Value = 210.25
